$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (Förändrad) rows 2 through 67 all hold the serial date value
# 45190 which needs to be updated to 45192, keeping existing formatting.
$ws.Range("C2:C67").Value = 45192
